{"js": "// The document has a run containing the misspelled text \"Fcrackkzip\"\n// (should read \"Fcrackzip\", i.e. the real tool name \"fcrackzip\"). The\n// edit fixes the typo by retyping at that spot, which also drags Word's\n// \"_GoBack\" last-edit bookmark from the end of the document to the\n// retyped location (splitting the run in two, with the bookmark start/\n// end pair sitting between \"Fcrack\" and \"zip\").\n\nconst body = context.document.body;\n\n// 1) Remove the old \"_GoBack\" bookmark (it currently sits in the empty\n//    paragraph at the very end of the document).\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// 2) Fix the typo: \"Fcrackkzip\" -> \"Fcrackzip\".\nconst misspelled = body.search(\"Fcrackkzip\", { matchCase: true });\nmisspelled.load(\"items\");\nawait context.sync();\nmisspelled.items[0].insertText(\"Fcrackzip\", \"Replace\");\nawait context.sync();\n\n// 3) Re-plant \"_GoBack\" right after \"Fcrack\", i.e. between \"Fcrack\" and\n//    \"zip\" - this is what splits the corrected word into two runs.\nconst fcrack = body.search(\"Fcrack\", { matchCase: true });\nfcrack.load(\"items\");\nawait context.sync();\nconst splitPoint = fcrack.items[0].getRange(\"End\");\nsplitPoint.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# The document has a run containing the misspelled text \"Fcrackkzip\"\n# (should read \"Fcrackzip\", i.e. the real tool name \"fcrackzip\"). The\n# edit fixes the typo by retyping at that spot, which also drags Word's\n# \"_GoBack\" last-edit bookmark from the end of the document to the\n# retyped location (splitting the run in two, with the bookmark start/\n# end pair sitting between \"Fcrack\" and \"zip\").\n\n$d = $word.ActiveDocument\n\n# 1) Remove the old \"_GoBack\" bookmark (it currently sits in the empty\n#    paragraph at the very end of the document).\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# 2) Fix the typo: \"Fcrackkzip\" -> \"Fcrackzip\".\n$rng = $d.Content\n$rng.Find.Text = \"Fcrackkzip\"\n$rng.Find.MatchCase = $true\n$rng.Find.Forward = $true\n$rng.Find.Execute() | Out-Null\n$rng.Text = \"Fcrackzip\"\n\n# 3) Re-plant \"_GoBack\" right after \"Fcrack\", i.e. between \"Fcrack\" and\n#    \"zip\" - this is what splits the corrected word into two runs.\n$rng2 = $d.Content\n$rng2.Find.Text = \"Fcrack\"\n$rng2.Find.MatchCase = $true\n$rng2.Find.Forward = $true\n$rng2.Find.Execute() | Out-Null\n$rng2.Collapse(0)  # wdCollapseEnd\n$d.Bookmarks.Add(\"_GoBack\", $rng2) | Out-Null\n"}
